# Auto-generated Excel COM-interop script applying the Sargatanas_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 6495.75  # H18: 5256.6 -> 6495.75
$ws.Cells.Item(18, 9).Value = 6862.8335  # I18: 5855.2856 -> 6862.8335
$ws.Cells.Item(18, 10).Value = 5394.5  # J18: 3859.6667 -> 5394.5
$ws.Cells.Item(18, 11).Value = 6862.8335  # K18: 5855.2856 -> 6862.8335
$ws.Cells.Item(18, 12).Value = 5394.5  # L18: 3859.6667 -> 5394.5
$ws.Cells.Item(18, 13).Value = -6578.8335  # M18: -5571.2856 -> -6578.8335
$ws.Cells.Item(18, 14).Value = -5962.5  # N18: -4427.6667 -> -5962.5

$ws.Cells.Item(33, 8).Value = 848.05884  # H33: 791.2222 -> 848.05884
$ws.Cells.Item(33, 9).Value = 963.0769  # I33: 852.93335 -> 963.0769
$ws.Cells.Item(33, 10).Value = 474.25  # J33: 482.66666 -> 474.25
$ws.Cells.Item(33, 11).Value = 963.0769  # K33: 852.93335 -> 963.0769
$ws.Cells.Item(33, 12).Value = 474.25  # L33: 482.66666 -> 474.25
$ws.Cells.Item(33, 13).Value = -734.0769  # M33: -623.93335 -> -734.0769
$ws.Cells.Item(33, 14).Value = -932.25  # N33: -940.66666 -> -932.25

$ws.Cells.Item(92, 8).Value = 3554.2727  # H92: 5056.857 -> 3554.2727
$ws.Cells.Item(92, 9).Value = 1140.8334  # I92: 1599 -> 1140.8334
$ws.Cells.Item(92, 10).Value = 6450.4  # J92: 9667.333000000001 -> 6450.4
$ws.Cells.Item(92, 11).Value = 1140.8334  # K92: 1599 -> 1140.8334
$ws.Cells.Item(92, 12).Value = 6450.4  # L92: 9667.333000000001 -> 6450.4
$ws.Cells.Item(92, 13).Value = 107.1666  # M92: -351 -> 107.1666
$ws.Cells.Item(92, 14).Value = -8946.4  # N92: -12163.333 -> -8946.4

$ws.Cells.Item(106, 8).Value = 3655  # H106: 3599.5 -> 3655
$ws.Cells.Item(106, 9).Value = 3655  # I106: 3599.5 -> 3655
$ws.Cells.Item(106, 11).Value = 3655  # K106: 3599.5 -> 3655
$ws.Cells.Item(106, 13).Value = -3024  # M106: -2968.5 -> -3024

$ws.Cells.Item(107, 8).Value = 46879450  # H107: 18751814 -> 46879450
$ws.Cells.Item(107, 9).Value = 17862228  # I107: 6946432 -> 17862228
$ws.Cells.Item(107, 10).Value = 250000000  # J107: 125000250 -> 250000000
$ws.Cells.Item(107, 11).Value = 17862228  # K107: 6946432 -> 17862228
$ws.Cells.Item(107, 12).Value = 250000000  # L107: 125000250 -> 250000000
$ws.Cells.Item(107, 13).Value = -17860308  # M107: -6944512 -> -17860308
$ws.Cells.Item(107, 14).Value = -250003840  # N107: -125004090 -> -250003840

$ws.Cells.Item(118, 8).Value = 1574.7  # H118: 1504.2727 -> 1574.7
$ws.Cells.Item(118, 9).Value = 697.4286  # I118: 710.25 -> 697.4286
$ws.Cells.Item(118, 11).Value = 2092.2858  # K118: 2130.75 -> 2092.2858
$ws.Cells.Item(118, 13).Value = -435.2857999999997  # M118: -473.75 -> -435.2857999999997

$ws.Cells.Item(129, 8).Value = 1300  # H129: 1238.2858 -> 1300
$ws.Cells.Item(129, 9).Value = 1000  # I129: 1194.6666 -> 1000
$ws.Cells.Item(129, 10).Value = 1750  # J129: 1500 -> 1750
$ws.Cells.Item(129, 11).Value = 3000  # K129: 3583.9998 -> 3000
$ws.Cells.Item(129, 12).Value = 5250  # L129: 4500 -> 5250
$ws.Cells.Item(129, 13).Value = 2000  # M129: 1416.0002 -> 2000
$ws.Cells.Item(129, 14).Value = -15250  # N129: -14500 -> -15250

$ws.Cells.Item(141, 8).Value = 2358.9565  # H141: 2172.1482 -> 2358.9565
$ws.Cells.Item(141, 9).Value = 1710.9333  # I141: 1581.8948 -> 1710.9333
$ws.Cells.Item(141, 11).Value = 5132.7999  # K141: 4745.6844 -> 5132.7999
$ws.Cells.Item(141, 13).Value = 47.20010000000002  # M141: 434.3155999999999 -> 47.20010000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2484.2964  # H102: 2652.4443 -> 2484.2964
$ws.Cells.Item(102, 9).Value = 2387.5386  # I102: 2456.64 -> 2387.5386
$ws.Cells.Item(102, 10).Value = 5000  # J102: 5100 -> 5000
$ws.Cells.Item(102, 11).Value = 2387.5386  # K102: 2456.64 -> 2387.5386
$ws.Cells.Item(102, 12).Value = 5000  # L102: 5100 -> 5000
$ws.Cells.Item(102, 13).Value = -765.5385999999999  # M102: -834.6399999999999 -> -765.5385999999999
$ws.Cells.Item(102, 14).Value = -8244  # N102: -8344 -> -8244

$ws.Cells.Item(122, 8).Value = 4428.222  # H122: 4331.25 -> 4428.222
$ws.Cells.Item(122, 9).Value = 3409.9  # I122: 3418 -> 3409.9
$ws.Cells.Item(122, 11).Value = 10229.7  # K122: 10254 -> 10229.7
$ws.Cells.Item(122, 13).Value = -7779.700000000001  # M122: -7804 -> -7779.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7937960  # H20: 7247746 -> 7937960
$ws.Cells.Item(20, 9).Value = 9805389  # I20: 8773295 -> 9805389
$ws.Cells.Item(20, 11).Value = 9805389  # K20: 8773295 -> 9805389
$ws.Cells.Item(20, 13).Value = -9805142  # M20: -8773048 -> -9805142

$ws.Cells.Item(99, 8).Value = 22730562  # H99: 30306012 -> 22730562
$ws.Cells.Item(99, 9).Value = 4330  # I99: 4450 -> 4330
$ws.Cells.Item(99, 11).Value = 4330  # K99: 4450 -> 4330
$ws.Cells.Item(99, 13).Value = -2832  # M99: -2952 -> -2832

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5065.7856  # H16: 4477.5625 -> 5065.7856
$ws.Cells.Item(16, 9).Value = 2272  # I16: 1794 -> 2272
$ws.Cells.Item(16, 11).Value = 2272  # K16: 1794 -> 2272
$ws.Cells.Item(16, 13).Value = -1985  # M16: -1507 -> -1985

$ws.Cells.Item(31, 8).Value = 10107.97  # H31: 10054.912 -> 10107.97
$ws.Cells.Item(31, 10).Value = 12286.913  # J31: 12120.958 -> 12286.913
$ws.Cells.Item(31, 12).Value = 12286.913  # L31: 12120.958 -> 12286.913
$ws.Cells.Item(31, 14).Value = -12876.913  # N31: -12710.958 -> -12876.913

$ws.Cells.Item(34, 8).Value = 10107.97  # H34: 10054.912 -> 10107.97
$ws.Cells.Item(34, 10).Value = 12286.913  # J34: 12120.958 -> 12286.913
$ws.Cells.Item(34, 12).Value = 12286.913  # L34: 12120.958 -> 12286.913
$ws.Cells.Item(34, 14).Value = -12690.913  # N34: -12524.958 -> -12690.913

$ws.Cells.Item(105, 8).Value = 5496360  # H105: 5954331.5 -> 5496360
$ws.Cells.Item(105, 9).Value = 8929659  # I105: 10205225 -> 8929659
$ws.Cells.Item(105, 11).Value = 8929659  # K105: 10205225 -> 8929659
$ws.Cells.Item(105, 13).Value = -8927912  # M105: -10203478 -> -8927912

$ws.Cells.Item(107, 8).Value = 1653.4166  # H107: 1660.4872 -> 1653.4166
$ws.Cells.Item(107, 9).Value = 1567.7142  # I107: 1631.5714 -> 1567.7142
$ws.Cells.Item(107, 10).Value = 1773.4  # J107: 1694.2222 -> 1773.4
$ws.Cells.Item(107, 11).Value = 1567.7142  # K107: 1631.5714 -> 1567.7142
$ws.Cells.Item(107, 12).Value = 1773.4  # L107: 1694.2222 -> 1773.4
$ws.Cells.Item(107, 13).Value = 352.2858000000001  # M107: 288.4286 -> 352.2858000000001
$ws.Cells.Item(107, 14).Value = -5613.4  # N107: -5534.2222 -> -5613.4

$ws.Cells.Item(113, 8).Value = 5065.7856  # H113: 4477.5625 -> 5065.7856
$ws.Cells.Item(113, 9).Value = 2272  # I113: 1794 -> 2272
$ws.Cells.Item(113, 11).Value = 2272  # K113: 1794 -> 2272
$ws.Cells.Item(113, 13).Value = -102  # M113: 376 -> -102

$ws.Cells.Item(132, 8).Value = 4724.2456  # H132: 4755.0537 -> 4724.2456
$ws.Cells.Item(132, 10).Value = 8126.6113  # J132: 8428.235000000001 -> 8126.6113
$ws.Cells.Item(132, 12).Value = 24379.8339  # L132: 25284.705 -> 24379.8339
$ws.Cells.Item(132, 14).Value = -29439.8339  # N132: -30344.705 -> -29439.8339

$ws.Cells.Item(138, 8).Value = 79999  # H138: 99999.5 -> 79999
$ws.Cells.Item(138, 10).Value = 79999  # J138: 99999.5 -> 79999
$ws.Cells.Item(138, 12).Value = 79999  # L138: 99999.5 -> 79999
$ws.Cells.Item(138, 14).Value = -90279  # N138: -110279.5 -> -90279

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 5651.4707  # H34: 5393.0557 -> 5651.4707
$ws.Cells.Item(34, 10).Value = 6847.4287  # J34: 6457.6 -> 6847.4287
$ws.Cells.Item(34, 12).Value = 20542.2861  # L34: 19372.8 -> 20542.2861
$ws.Cells.Item(34, 14).Value = -20710.2861  # N34: -19540.8 -> -20710.2861

$ws.Cells.Item(93, 8).Value = 6922.4  # H93: 7136 -> 6922.4
$ws.Cells.Item(93, 10).Value = 6922.222  # J93: 7162.5 -> 6922.222
$ws.Cells.Item(93, 12).Value = 20766.666  # L93: 21487.5 -> 20766.666
$ws.Cells.Item(93, 14).Value = -24510.666  # N93: -25231.5 -> -24510.666

$ws.Cells.Item(97, 8).Value = 234.66667  # H97: 251 -> 234.66667
$ws.Cells.Item(97, 10).Value = 234.66667  # J97: 251 -> 234.66667
$ws.Cells.Item(97, 12).Value = 704.00001  # L97: 753 -> 704.00001
$ws.Cells.Item(97, 14).Value = -1696.00001  # N97: -1745 -> -1696.00001

$ws.Cells.Item(122, 8).Value = 2573483.5  # H122: 2830737.5 -> 2573483.5
$ws.Cells.Item(122, 10).Value = 2590  # J122: 2919.6 -> 2590
$ws.Cells.Item(122, 12).Value = 23310  # L122: 26276.4 -> 23310
$ws.Cells.Item(122, 14).Value = -28210  # N122: -31176.4 -> -28210

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8624.789000000001  # H70: 8624.842000000001 -> 8624.789000000001
$ws.Cells.Item(70, 9).Value = 7881.0835  # I70: 7881.1665 -> 7881.0835
$ws.Cells.Item(70, 11).Value = 7881.0835  # K70: 7881.1665 -> 7881.0835
$ws.Cells.Item(70, 13).Value = -7611.0835  # M70: -7611.1665 -> -7611.0835

$ws.Cells.Item(73, 8).Value = 8624.789000000001  # H73: 8624.842000000001 -> 8624.789000000001
$ws.Cells.Item(73, 9).Value = 7881.0835  # I73: 7881.1665 -> 7881.0835
$ws.Cells.Item(73, 11).Value = 7881.0835  # K73: 7881.1665 -> 7881.0835
$ws.Cells.Item(73, 13).Value = -6945.0835  # M73: -6945.1665 -> -6945.0835

$ws.Cells.Item(97, 8).Value = 1891.6  # H97: 1825.3334 -> 1891.6
$ws.Cells.Item(97, 10).Value = 1821.3334  # J97: 1689.2 -> 1821.3334
$ws.Cells.Item(97, 12).Value = 1821.3334  # L97: 1689.2 -> 1821.3334
$ws.Cells.Item(97, 14).Value = -2813.3334  # N97: -2681.2 -> -2813.3334

$ws.Cells.Item(113, 8).Value = 8631.772000000001  # H113: 6689.3667 -> 8631.772000000001
$ws.Cells.Item(113, 9).Value = 5979.8  # I113: 3588.7 -> 5979.8
$ws.Cells.Item(113, 10).Value = 9411.764999999999  # J113: 8239.700000000001 -> 9411.764999999999
$ws.Cells.Item(113, 11).Value = 5979.8  # K113: 3588.7 -> 5979.8
$ws.Cells.Item(113, 12).Value = 9411.764999999999  # L113: 8239.700000000001 -> 9411.764999999999
$ws.Cells.Item(113, 13).Value = -3809.8  # M113: -1418.7 -> -3809.8
$ws.Cells.Item(113, 14).Value = -13751.765  # N113: -12579.7 -> -13751.765

$ws.Cells.Item(122, 8).Value = 338833  # H122: 505999.5 -> 338833
$ws.Cells.Item(122, 10).Value = 4500  # J122: 0 -> 4500
$ws.Cells.Item(122, 12).Value = 13500  # L122: 0 -> 13500
$ws.Cells.Item(122, 14).Value = -18400  # N122: add

$ws.Cells.Item(132, 8).Value = 5870  # H132: 5551.625 -> 5870
$ws.Cells.Item(132, 9).Value = 1700  # I132: 1597.3334 -> 1700
$ws.Cells.Item(132, 11).Value = 5100  # K132: 4792.0002 -> 5100
$ws.Cells.Item(132, 13).Value = -2570  # M132: -2262.0002 -> -2570

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4489.25  # H7: 6493.1665 -> 4489.25
$ws.Cells.Item(7, 9).Value = 4489.25  # I7: 4979.5 -> 4489.25
$ws.Cells.Item(7, 10).Value = 0  # J7: 7250 -> 0
$ws.Cells.Item(7, 11).Value = 4489.25  # K7: 4979.5 -> 4489.25
$ws.Cells.Item(7, 12).Value = 0  # L7: 7250 -> 0
$ws.Cells.Item(7, 13).Value = -4377.25  # M7: -4867.5 -> -4377.25
$ws.Cells.Item(7, 14).ClearContents()  # N7: delete

$ws.Cells.Item(40, 8).Value = 4741.421  # H40: 4843.778 -> 4741.421
$ws.Cells.Item(40, 9).Value = 4005.8572  # I40: 4091 -> 4005.8572
$ws.Cells.Item(40, 11).Value = 4005.8572  # K40: 4091 -> 4005.8572
$ws.Cells.Item(40, 13).Value = -3869.8572  # M40: -3955 -> -3869.8572

$ws.Cells.Item(55, 8).Value = 55555900  # H55: 58823892 -> 55555900
$ws.Cells.Item(55, 9).Value = 100000080  # I55: 111111190 -> 100000080
$ws.Cells.Item(55, 11).Value = 100000080  # K55: 111111190 -> 100000080
$ws.Cells.Item(55, 13).Value = -99999907  # M55: -111111017 -> -99999907

$ws.Cells.Item(61, 8).Value = 3570.1292  # H61: 3669.2 -> 3570.1292
$ws.Cells.Item(61, 9).Value = 2450.0588  # I61: 2565.8125 -> 2450.0588
$ws.Cells.Item(61, 11).Value = 2450.0588  # K61: 2565.8125 -> 2450.0588
$ws.Cells.Item(61, 13).Value = -2248.0588  # M61: -2363.8125 -> -2248.0588

$ws.Cells.Item(113, 8).Value = 3570.1292  # H113: 3669.2 -> 3570.1292
$ws.Cells.Item(113, 9).Value = 2450.0588  # I113: 2565.8125 -> 2450.0588
$ws.Cells.Item(113, 11).Value = 2450.0588  # K113: 2565.8125 -> 2450.0588
$ws.Cells.Item(113, 13).Value = -280.0587999999998  # M113: -395.8125 -> -280.0587999999998

$ws.Cells.Item(126, 8).Value = 4489.25  # H126: 6493.1665 -> 4489.25
$ws.Cells.Item(126, 9).Value = 4489.25  # I126: 4979.5 -> 4489.25
$ws.Cells.Item(126, 10).Value = 0  # J126: 7250 -> 0
$ws.Cells.Item(126, 11).Value = 13467.75  # K126: 14938.5 -> 13467.75
$ws.Cells.Item(126, 12).Value = 0  # L126: 21750 -> 0
$ws.Cells.Item(126, 13).Value = -10997.75  # M126: -12468.5 -> -10997.75
$ws.Cells.Item(126, 14).ClearContents()  # N126: delete

$ws.Cells.Item(132, 8).Value = 16674568  # H132: 17249518 -> 16674568
$ws.Cells.Item(132, 9).Value = 45458044  # I132: 55559056 -> 45458044
$ws.Cells.Item(132, 10).Value = 10450.211  # J132: 10227.65 -> 10450.211
$ws.Cells.Item(132, 11).Value = 136374132  # K132: 166677168 -> 136374132
$ws.Cells.Item(132, 12).Value = 31350.633  # L132: 30682.95 -> 31350.633
$ws.Cells.Item(132, 13).Value = -136371602  # M132: -166674638 -> -136371602
$ws.Cells.Item(132, 14).Value = -36410.633  # N132: -35742.95 -> -36410.633

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 0  # H6: 4999 -> 0
$ws.Cells.Item(6, 10).Value = 0  # J6: 4999 -> 0
$ws.Cells.Item(6, 12).Value = 0  # L6: 4999 -> 0
$ws.Cells.Item(6, 14).ClearContents()  # N6: delete
